$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove worker "1032461151 / PAOLA CATALINA YEPES ARIZA" (6 rows, periods 2502-2507).
# This shifts the remaining worker (1143392460 / ANDREA CAROLINA ACEVEDO GUERRA, rows
# 23-26) up to rows 17-20, and the trailing signature block (old rows 31-32) up to 25-26,
# matching the new table layout.
$ws.Range("A17:A22").EntireRow.Delete()

# Update the (now only remaining) worker's periods to run 2107..2110 in ascending order
# and bump the "Salario Basico" (G column) for each period from 1035000 to 1486000.
$ws.Range("E17").Value = "2107"
$ws.Range("F17").Value = 41400
$ws.Range("G17").Value = 1486000

$ws.Range("E18").Value = "2108"
$ws.Range("F18").Value = 41400
$ws.Range("G18").Value = 1486000

$ws.Range("E19").Value = "2109"
$ws.Range("F19").Value = 41400
$ws.Range("G19").Value = 1486000

$ws.Range("E20").Value = "2110"
$ws.Range("F20").Value = 41400
$ws.Range("G20").Value = 1486000

# Refresh the summary figures: total "Valor Mora", worker count, and period count.
$ws.Range("E11").Value = 176267
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 5
